# Generate Report for Handback
# Row 7 in both the "zh-cn" and "de-de" sheets moves from "awaiting handback"
# to "handback received, but stale version" - fills in the target file,
# handback file, handback datetime and an error detail message, and adds a
# hyperlink on the newly-populated "Latest Target File" cell (column I).

$wb = $excel.ActiveWorkbook

$errorMessage = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/238432823cbb67767cc92e2b1c67437c9bb94948/e2e/3b49ae0b-3476-41e6-a471-094ecce6613b.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/67ba02d2369b9374a0c64675a89577ef3f81db9b/e2e/3b49ae0b-3476-41e6-a471-094ecce6613b.md."

# ---------------------------------------------------------------------
# zh-cn sheet, row 7
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("I7").Value = "3b49ae0b-3476-41e6-a471-094ecce6613b.md"
$wsZh.Range("J7").Value = "3b49ae0b-3476-41e6-a471-094ecce6613b.d77d4aa513d38468a8793962c6e3ec50e6743a51.zh-cn.xlf"
$wsZh.Range("K7").Value = "2016-09-07 11:17:04"
$wsZh.Range("P7").Value = $errorMessage

$wsZh.Hyperlinks.Add($wsZh.Range("I7"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/238432823cbb67767cc92e2b1c67437c9bb94948/e2e/3b49ae0b-3476-41e6-a471-094ecce6613b.md", "", "", "3b49ae0b-3476-41e6-a471-094ecce6613b.md")

# ---------------------------------------------------------------------
# de-de sheet, row 7
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("I7").Value = "3b49ae0b-3476-41e6-a471-094ecce6613b.md"
$wsDe.Range("J7").Value = "3b49ae0b-3476-41e6-a471-094ecce6613b.d77d4aa513d38468a8793962c6e3ec50e6743a51.de-de.xlf"
$wsDe.Range("K7").Value = "2016-09-07 11:17:19"
$wsDe.Range("P7").Value = $errorMessage

$wsDe.Hyperlinks.Add($wsDe.Range("I7"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/238432823cbb67767cc92e2b1c67437c9bb94948/e2e/3b49ae0b-3476-41e6-a471-094ecce6613b.md", "", "", "3b49ae0b-3476-41e6-a471-094ecce6613b.md")
